$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in "Tom" for the remaining use-case rows (matching rows 2-4 that already have it)
$ws.Range("B6").Value = "Tom"
$ws.Range("B8").Value = "Tom"
$ws.Range("B9").Value = "Tom"
$ws.Range("B10").Value = "Tom"

# Move the active selection to D9
$ws.Range("D9").Select()
